# Apply the updated cryptocurrency price/volume snapshot to the sheet.
# Source values come from the latest coinranking.com pull; cells keep their
# original 'text' formatting so values like 1.000 and 27.093.54 aren't
# reinterpreted as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = 'Normal'
}

# Row 2
Set-TextValue 'D2' '27.093.54'
Set-TextValue 'E2' '  -2.55%  '

# Row 3
Set-TextValue 'D3' '1.864.55'
Set-TextValue 'E3' '  -2.37%  '

# Row 4
Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  -0.05%  '

# Row 5
Set-TextValue 'D5' '306.15'
Set-TextValue 'E5' '  -2.10%  '

# Row 6
Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  +0.02%  '

# Row 7
Set-TextValue 'D7' '0.5160'
Set-TextValue 'E7' '  -0.24%  '

# Row 8
Set-TextValue 'D8' '0.3760'
Set-TextValue 'E8' '  -0.59%  '

# Row 9
Set-TextValue 'D9' '0.07154'
Set-TextValue 'E9' '  -1.38%  '

# Row 10
Set-TextValue 'D10' '0.8887'
Set-TextValue 'E10' '  -1.85%  '

# Row 11
Set-TextValue 'D11' '20.69'
Set-TextValue 'E11' '  -2.79%  '

# Row 12
Set-TextValue 'D12' '0.07550'
Set-TextValue 'E12' '  -1.41%  '

# Row 13
Set-TextValue 'D13' '1.866.27'
Set-TextValue 'E13' '  -2.60%  '

# Row 14
Set-TextValue 'E14' '  -2.73%  '

# Row 15
Set-TextValue 'D15' '89.61'
Set-TextValue 'E15' '  -2.80%  '

# Row 16
Set-TextValue 'D16' '1.001'
Set-TextValue 'E16' '  +0.01%  '

# Row 17
Set-TextValue 'D17' '0.000008470'
Set-TextValue 'E17' '  -2.77%  '

# Row 18
Set-TextValue 'D18' '14.07'
Set-TextValue 'E18' '  -3.27%  '

# Row 19
Set-TextValue 'D19' '1.001'
Set-TextValue 'E19' '  +0.05%  '

# Row 20
Set-TextValue 'D20' '27.120.88'
Set-TextValue 'E20' '  -2.61%  '

# Row 21
Set-TextValue 'D21' '5.009'
Set-TextValue 'E21' '  -2.83%  '

# Row 22
Set-TextValue 'D22' '2.079.16'
Set-TextValue 'E22' '  -4.60%  '

# Row 23
Set-TextValue 'E23' '  -3.57%  '

# Row 24
Set-TextValue 'D24' '6.443'
Set-TextValue 'E24' '  -3.02%  '

# Row 25
Set-TextValue 'D25' '1.837'
Set-TextValue 'E25' '  -1.85%  '

# Row 26
Set-TextValue 'D26' '145.22'
Set-TextValue 'E26' '  -5.74%  '

# Row 27
Set-TextValue 'E27' '  -2.29%  '

# Row 28
Set-TextValue 'D28' '2.086'
Set-TextValue 'E28' '  -4.00%  '

# Row 29
Set-TextValue 'D29' '112.74'
Set-TextValue 'E29' '  -2.00%  '

# Row 30
Set-TextValue 'D30' '4.659'
Set-TextValue 'E30' '  -4.21%  '

# Row 31
Set-TextValue 'D31' '4.662'
Set-TextValue 'E31' '  -3.88%  '

# Row 32
Set-TextValue 'D32' '0.09162'
Set-TextValue 'E32' '  +0.89%  '

# Row 33
Set-TextValue 'D33' '0.05099'
Set-TextValue 'E33' '  -3.70%  '

# Row 34
Set-TextValue 'D34' '3.072'
Set-TextValue 'E34' '  -3.59%  '

# Row 35
Set-TextValue 'E35' '  -6.44%  '

# Row 36
Set-TextValue 'D36' '0.7234'
Set-TextValue 'E36' '  -7.40%  '

# Row 37
Set-TextValue 'D37' '0.02037'
Set-TextValue 'E37' '  -2.89%  '

# Row 38
Set-TextValue 'D38' '3.084'
Set-TextValue 'E38' '  +0.16%  '

# Row 39
Set-TextValue 'D39' '2.486'
Set-TextValue 'E39' '  -4.83%  '

# Row 40
Set-TextValue 'E40' '  -1.77%  '

# Row 41
Set-TextValue 'D41' '0.5272'
Set-TextValue 'E41' '  -5.68%  '

# Row 42
Set-TextValue 'B42' 'Quant'
Set-TextValue 'C42' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D42' '116.25'
Set-TextValue 'E42' '  +0.75%  '

# Row 43
Set-TextValue 'B43' 'FraxShare'
Set-TextValue 'C43' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D43' '6.460'
Set-TextValue 'E43' '  -4.04%  '

# Row 44
Set-TextValue 'D44' '8.276'
Set-TextValue 'E44' '  -3.42%  '

# Row 45
Set-TextValue 'D45' '0.1464'
Set-TextValue 'E45' '  -3.57%  '

# Row 46
Set-TextValue 'E46' '  +0.05%  '

# Row 48
Set-TextValue 'D48' '9.938'
Set-TextValue 'E48' '  -5.25%  '

# Row 49
Set-TextValue 'D49' '1.563'
Set-TextValue 'E49' '  -3.56%  '

# Row 50
Set-TextValue 'E50' '  -1.41%  '

# Row 51
Set-TextValue 'E51' '  -5.51%  '
